$wb = $excel.ActiveWorkbook

# --- Update the raw metric values on the "Metrics" sheet ---
$wsMetrics = $wb.Worksheets.Item("Metrics")

$wsMetrics.Range("B2").Value = 259787.27
$wsMetrics.Range("B3").Value = 213263.62000000002
$wsMetrics.Range("B4").Value = 82586.849999999991
$wsMetrics.Range("B5").Value = 10367
$wsMetrics.Range("B6").Value = 4626918.74
$wsMetrics.Range("B7").Value = 3903082.2899999996
$wsMetrics.Range("B8").Value = 1353188.99
$wsMetrics.Range("B9").Value = 179368
$wsMetrics.Range("B10").Value = 33092242.540999822
$wsMetrics.Range("B11").Value = 31178303.810000002
$wsMetrics.Range("B12").Value = 11634897.880000001
$wsMetrics.Range("B13").Value = 1276995

# Move the cached selection on "Metrics" to match the latest edit location.
$wsMetrics.Activate() | Out-Null
$wsMetrics.Range("E17").Select() | Out-Null

# --- "today" sheet formulas (Metrics!B2..B13 references, plus derived E/F
# columns) recompute automatically on recalc. Just restore it as the active
# sheet and update its cached selection. ---
$wsToday = $wb.Worksheets.Item("today")
$wsToday.Activate() | Out-Null
$wsToday.Range("H1").Select() | Out-Null
